# Auto-generated edit script applying the Durandal_Profits.xlsx diff
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H-N) across 8 sheets
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1107.5
$ws.Range("I40").Value = 1066.1072
$ws.Range("J40").Value = 1204.0834
$ws.Range("K40").Value = 1066.1072
$ws.Range("L40").Value = 1204.0834
$ws.Range("M40").Value = -891.1071999999999
$ws.Range("N40").Value = -1554.0834
$ws.Range("H129").Value = 862.48
$ws.Range("I129").Value = 288.58334
$ws.Range("K129").Value = 865.7500200000001
$ws.Range("M129").Value = 4134.24998

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1101.56
$ws.Range("I61").Value = 870.6875
$ws.Range("J61").Value = 1512
$ws.Range("K61").Value = 870.6875
$ws.Range("L61").Value = 1512
$ws.Range("M61").Value = -658.6875
$ws.Range("N61").Value = -1936
$ws.Range("H74").Value = 2003.5
$ws.Range("I74").Value = 1617
$ws.Range("J74").Value = 2969.75
$ws.Range("K74").Value = 1617
$ws.Range("L74").Value = 2969.75
$ws.Range("M74").Value = -743
$ws.Range("N74").Value = -4717.75
$ws.Range("H77").Value = 2003.5
$ws.Range("I77").Value = 1617
$ws.Range("J77").Value = 2969.75
$ws.Range("K77").Value = 8085
$ws.Range("L77").Value = 14848.75
$ws.Range("M77").Value = -3717
$ws.Range("N77").Value = -23584.75
$ws.Range("H102").Value = 3213.28
$ws.Range("I102").Value = 2491.476
$ws.Range("K102").Value = 2491.476
$ws.Range("M102").Value = -869.4760000000001
$ws.Range("H132").Value = 1066.6957
$ws.Range("I132").Value = 767
$ws.Range("J132").Value = 2490.25
$ws.Range("K132").Value = 2301
$ws.Range("L132").Value = 7470.75
$ws.Range("M132").Value = 229
$ws.Range("N132").Value = -12530.75
$ws.Range("H136").Value = 1101.56
$ws.Range("I136").Value = 870.6875
$ws.Range("J136").Value = 1512
$ws.Range("K136").Value = 2612.0625
$ws.Range("L136").Value = 4536
$ws.Range("M136").Value = -62.0625
$ws.Range("N136").Value = -9636

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 2059.2222
$ws.Range("I80").Value = 310.8
$ws.Range("J80").Value = 2731.6924
$ws.Range("K80").Value = 310.8
$ws.Range("L80").Value = 2731.6924
$ws.Range("M80").Value = 687.2
$ws.Range("N80").Value = -4727.6924
$ws.Range("H83").Value = 2059.2222
$ws.Range("I83").Value = 310.8
$ws.Range("J83").Value = 2731.6924
$ws.Range("K83").Value = 1554
$ws.Range("L83").Value = 13658.462
$ws.Range("M83").Value = 3438
$ws.Range("N83").Value = -23642.462
$ws.Range("H134").Value = 7067.6665
$ws.Range("I134").Value = 838.1579
$ws.Range("K134").Value = 2514.4737
$ws.Range("M134").Value = 20.52629999999999

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8476593
$ws.Range("I31").Value = 11629495
$ws.Range("J31").Value = 3171
$ws.Range("K31").Value = 11629495
$ws.Range("L31").Value = 3171
$ws.Range("M31").Value = -11629200
$ws.Range("N31").Value = -3761
$ws.Range("H34").Value = 8476593
$ws.Range("I34").Value = 11629495
$ws.Range("J34").Value = 3171
$ws.Range("K34").Value = 11629495
$ws.Range("L34").Value = 3171
$ws.Range("M34").Value = -11629293
$ws.Range("N34").Value = -3575
$ws.Range("H58").Value = 1200.9
$ws.Range("I58").Value = 940.05
$ws.Range("K58").Value = 940.05
$ws.Range("M58").Value = -737.05
$ws.Range("H136").Value = 1200.9
$ws.Range("I136").Value = 940.05
$ws.Range("K136").Value = 2820.15
$ws.Range("M136").Value = -270.1499999999996

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 10870577
$ws.Range("I131").Value = 677.61536
$ws.Range("J131").Value = 15152658
$ws.Range("K131").Value = 2032.84608
$ws.Range("L131").Value = 45457974
$ws.Range("M131").Value = 3007.15392
$ws.Range("N131").Value = -45468054
$ws.Range("H136").Value = 2000
$ws.Range("I136").Value = 2000
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 6000
$ws.Range("L136").Value = 0
$ws.Range("M136").ClearContents()
$ws.Range("N136").Value = -900

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 64417.945
$ws.Range("I80").Value = 224401
$ws.Range("J80").Value = 2886
$ws.Range("K80").Value = 224401
$ws.Range("L80").Value = 2886
$ws.Range("M80").Value = -223403
$ws.Range("N80").Value = -4882
$ws.Range("H83").Value = 64417.945
$ws.Range("I83").Value = 224401
$ws.Range("J83").Value = 2886
$ws.Range("K83").Value = 1122005
$ws.Range("L83").Value = 14430
$ws.Range("M83").Value = -1117013
$ws.Range("N83").Value = -24414
$ws.Range("H132").Value = 36430.137
$ws.Range("I132").Value = 38906.445
$ws.Range("K132").Value = 116719.335
$ws.Range("M132").Value = -114189.335
$ws.Range("H133").Value = 48775.555
$ws.Range("J133").Value = 48775.555
$ws.Range("L133").Value = 48775.555
$ws.Range("N133").Value = -58895.555

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1562.4231
$ws.Range("I61").Value = 1315.1818
$ws.Range("J61").Value = 1743.7333
$ws.Range("K61").Value = 1315.1818
$ws.Range("L61").Value = 1743.7333
$ws.Range("M61").Value = -1113.1818
$ws.Range("N61").Value = -2147.7333
$ws.Range("H68").Value = 2074.111
$ws.Range("I68").Value = 1993.3334
$ws.Range("K68").Value = 1993.3334
$ws.Range("M68").Value = -1244.3334
$ws.Range("H71").Value = 2074.111
$ws.Range("I71").Value = 1993.3334
$ws.Range("K71").Value = 9966.666999999999
$ws.Range("M71").Value = -6222.666999999999
$ws.Range("H113").Value = 1562.4231
$ws.Range("I113").Value = 1315.1818
$ws.Range("J113").Value = 1743.7333
$ws.Range("K113").Value = 1315.1818
$ws.Range("L113").Value = 1743.7333
$ws.Range("M113").Value = 854.8181999999999
$ws.Range("N113").Value = -6083.7333
$ws.Range("H130").Value = 33999.5
$ws.Range("J130").Value = 33999.5
$ws.Range("L130").Value = 33999.5
$ws.Range("N130").Value = -44039.5
$ws.Range("H136").Value = 2090.6038
$ws.Range("I136").Value = 1363.2903
$ws.Range("J136").Value = 3115.4546
$ws.Range("K136").Value = 4089.8709
$ws.Range("L136").Value = 9346.363799999999
$ws.Range("M136").Value = -1539.8709
$ws.Range("N136").Value = -14446.3638

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 8555.143
$ws.Range("I62").Value = 4064.182
$ws.Range("J62").Value = 13495.2
$ws.Range("K62").Value = 4064.182
$ws.Range("L62").Value = 13495.2
$ws.Range("M62").Value = -3440.182
$ws.Range("N62").Value = -14743.2
$ws.Range("H65").Value = 8555.143
$ws.Range("I65").Value = 4064.182
$ws.Range("J65").Value = 13495.2
$ws.Range("K65").Value = 20320.91
$ws.Range("L65").Value = 67476
$ws.Range("M65").Value = -17200.91
$ws.Range("N65").Value = -73716
$ws.Range("H103").Value = 29483
$ws.Range("J103").Value = 29483
$ws.Range("L103").Value = 29483
$ws.Range("N103").Value = -31827
$ws.Range("H132").Value = 15434216
$ws.Range("I132").Value = 21187338
$ws.Range("J132").Value = 5390.4546
$ws.Range("K132").Value = 63562014
$ws.Range("L132").Value = 16171.3638
$ws.Range("M132").Value = -63559484
$ws.Range("N132").Value = -21231.3638
$ws.Range("H136").Value = 2608.8542
$ws.Range("I136").Value = 3404.7878
$ws.Range("J136").Value = 857.8
$ws.Range("K136").Value = 10214.3634
$ws.Range("L136").Value = 2573.4
$ws.Range("M136").Value = -7664.3634
$ws.Range("N136").Value = -7673.4
